$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17,8).Value = 1458.6307
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 1458.6307
$ws.Cells.Item(17,11).Value = 0
$ws.Cells.Item(17,12).Value = 4375.8921
$ws.Cells.Item(17,14).Value = -4711.8921

# Row 21
$ws.Cells.Item(21,8).Value = 30506.334
$ws.Cells.Item(21,9).Value = 80019
$ws.Cells.Item(21,10).Value = 20603.8
$ws.Cells.Item(21,11).Value = 80019
$ws.Cells.Item(21,12).Value = 20603.8
$ws.Cells.Item(21,13).Value = -79551
$ws.Cells.Item(21,14).Value = -21539.8

# Row 23
$ws.Cells.Item(23,8).Value = 30506.334
$ws.Cells.Item(23,9).Value = 80019
$ws.Cells.Item(23,10).Value = 20603.8
$ws.Cells.Item(23,11).Value = 80019
$ws.Cells.Item(23,12).Value = 20603.8
$ws.Cells.Item(23,13).Value = -79785
$ws.Cells.Item(23,14).Value = -21071.8

# Row 49
$ws.Cells.Item(49,8).Value = 3500
$ws.Cells.Item(49,9).Value = 3500
$ws.Cells.Item(49,10).Value = 0
$ws.Cells.Item(49,11).Value = 10500
$ws.Cells.Item(49,12).Value = 0
$ws.Cells.Item(49,13).Value = -10364
$ws.Cells.Item(49,14).ClearContents()

# Row 76
$ws.Cells.Item(76,8).Value = 2930.7273
$ws.Cells.Item(76,9).Value = 2803.7368
$ws.Cells.Item(76,10).Value = 3735
$ws.Cells.Item(76,11).Value = 2803.7368
$ws.Cells.Item(76,12).Value = 3735
$ws.Cells.Item(76,13).Value = -2488.7368

# Row 79
$ws.Cells.Item(79,8).Value = 2930.7273
$ws.Cells.Item(79,9).Value = 2803.7368
$ws.Cells.Item(79,10).Value = 3735
$ws.Cells.Item(79,11).Value = 2803.7368
$ws.Cells.Item(79,12).Value = 3735
$ws.Cells.Item(79,13).Value = -1711.7368

# Row 107
$ws.Cells.Item(107,8).Value = 1457.5
$ws.Cells.Item(107,9).Value = 1370.5
$ws.Cells.Item(107,10).Value = 1675
$ws.Cells.Item(107,11).Value = 1370.5
$ws.Cells.Item(107,12).Value = 1675
$ws.Cells.Item(107,13).Value = 549.5

# Row 132
$ws.Cells.Item(132,8).Value = 3775904.8
$ws.Cells.Item(132,9).Value = 4446558.5
$ws.Cells.Item(132,10).Value = 3476.25
$ws.Cells.Item(132,11).Value = 13339675.5
$ws.Cells.Item(132,12).Value = 10428.75
$ws.Cells.Item(132,13).Value = -13337145.5
$ws.Cells.Item(132,14).Value = -15488.75

# Row 133
$ws.Cells.Item(133,8).Value = 29738.182
$ws.Cells.Item(133,9).Value = 0
$ws.Cells.Item(133,10).Value = 29738.182
$ws.Cells.Item(133,11).Value = 0
$ws.Cells.Item(133,12).Value = 29738.182
$ws.Cells.Item(133,14).Value = -39858.182

# Row 137
$ws.Cells.Item(137,8).Value = 3471.9375
$ws.Cells.Item(137,9).Value = 3517.25
$ws.Cells.Item(137,10).Value = 3336
$ws.Cells.Item(137,11).Value = 10551.75
$ws.Cells.Item(137,12).Value = 10008
$ws.Cells.Item(137,13).Value = -8001.75
$ws.Cells.Item(137,14).Value = -15108

# Row 138
$ws.Cells.Item(138,8).Value = 2379.6262
$ws.Cells.Item(138,9).Value = 929.5484
$ws.Cells.Item(138,10).Value = 3040.6912
$ws.Cells.Item(138,11).Value = 2788.6452
$ws.Cells.Item(138,12).Value = 9122.0736
$ws.Cells.Item(138,13).Value = 2351.3548
$ws.Cells.Item(138,14).Value = -19402.0736

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32,8).Value = 1441.67
$ws.Cells.Item(32,9).Value = 1282.1976
$ws.Cells.Item(32,10).Value = 2421.2856
$ws.Cells.Item(32,11).Value = 1282.1976
$ws.Cells.Item(32,12).Value = 2421.2856
$ws.Cells.Item(32,13).Value = -995.1976
$ws.Cells.Item(32,14).Value = -2995.2856

# Row 61
$ws.Cells.Item(61,8).Value = 2879.875
$ws.Cells.Item(61,9).Value = 1136.5555
$ws.Cells.Item(61,10).Value = 5121.2856
$ws.Cells.Item(61,11).Value = 1136.5555
$ws.Cells.Item(61,12).Value = 5121.2856
$ws.Cells.Item(61,13).Value = -924.5554999999999
$ws.Cells.Item(61,14).Value = -5545.2856

# Row 74
$ws.Cells.Item(74,8).Value = 626.2963
$ws.Cells.Item(74,9).Value = 587.0833
$ws.Cells.Item(74,10).Value = 940
$ws.Cells.Item(74,11).Value = 587.0833
$ws.Cells.Item(74,12).Value = 940
$ws.Cells.Item(74,13).Value = 286.9167
$ws.Cells.Item(74,14).Value = -2688

# Row 77
$ws.Cells.Item(77,8).Value = 626.2963
$ws.Cells.Item(77,9).Value = 587.0833
$ws.Cells.Item(77,10).Value = 940
$ws.Cells.Item(77,11).Value = 2935.4165
$ws.Cells.Item(77,12).Value = 4700
$ws.Cells.Item(77,13).Value = 1432.5835
$ws.Cells.Item(77,14).Value = -13436

# Row 86
$ws.Cells.Item(86,8).Value = 5285
$ws.Cells.Item(86,9).Value = 5285
$ws.Cells.Item(86,10).Value = 0
$ws.Cells.Item(86,11).Value = 5285
$ws.Cells.Item(86,12).Value = 0
$ws.Cells.Item(86,13).Value = -4099
$ws.Cells.Item(86,14).ClearContents()

# Row 89
$ws.Cells.Item(89,8).Value = 5285
$ws.Cells.Item(89,9).Value = 5285
$ws.Cells.Item(89,10).Value = 0
$ws.Cells.Item(89,11).Value = 15855
$ws.Cells.Item(89,12).Value = 0
$ws.Cells.Item(89,13).Value = -9927
$ws.Cells.Item(89,14).ClearContents()

# Row 132
$ws.Cells.Item(132,8).Value = 1718.4154
$ws.Cells.Item(132,9).Value = 1128.3877
$ws.Cells.Item(132,10).Value = 3525.375
$ws.Cells.Item(132,11).Value = 3385.1631
$ws.Cells.Item(132,12).Value = 10576.125
$ws.Cells.Item(132,13).Value = -855.1630999999998
$ws.Cells.Item(132,14).Value = -15636.125

# Row 136
$ws.Cells.Item(136,8).Value = 2879.875
$ws.Cells.Item(136,9).Value = 1136.5555
$ws.Cells.Item(136,10).Value = 5121.2856
$ws.Cells.Item(136,11).Value = 3409.6665
$ws.Cells.Item(136,12).Value = 15363.8568
$ws.Cells.Item(136,13).Value = -859.6664999999998
$ws.Cells.Item(136,14).Value = -20463.8568

# Row 139
$ws.Cells.Item(139,8).Value = 27500
$ws.Cells.Item(139,9).Value = 0
$ws.Cells.Item(139,10).Value = 27500
$ws.Cells.Item(139,11).Value = 0
$ws.Cells.Item(139,12).Value = 27500
$ws.Cells.Item(139,14).Value = -37780

# Row 141
$ws.Cells.Item(141,8).Value = 0
$ws.Cells.Item(141,9).Value = 0
$ws.Cells.Item(141,10).Value = 0
$ws.Cells.Item(141,11).Value = 0
$ws.Cells.Item(141,12).Value = 0
$ws.Cells.Item(141,14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86,8).Value = 556811.75
$ws.Cells.Item(86,9).Value = 959275.7
$ws.Cells.Item(86,10).Value = 39358.07
$ws.Cells.Item(86,11).Value = 959275.7
$ws.Cells.Item(86,12).Value = 39358.07
$ws.Cells.Item(86,13).Value = -958152.7
$ws.Cells.Item(86,14).Value = -41604.07

# Row 89
$ws.Cells.Item(89,8).Value = 556811.75
$ws.Cells.Item(89,9).Value = 959275.7
$ws.Cells.Item(89,10).Value = 39358.07
$ws.Cells.Item(89,11).Value = 4796378.5
$ws.Cells.Item(89,12).Value = 196790.35
$ws.Cells.Item(89,13).Value = -4790762.5
$ws.Cells.Item(89,14).Value = -208022.35

# Row 134
$ws.Cells.Item(134,8).Value = 2335.8667
$ws.Cells.Item(134,9).Value = 1726.8096
$ws.Cells.Item(134,10).Value = 3757
$ws.Cells.Item(134,11).Value = 5180.4288
$ws.Cells.Item(134,12).Value = 11271
$ws.Cells.Item(134,13).Value = -2645.4288
$ws.Cells.Item(134,14).Value = -16341

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31,8).Value = 3103.558
$ws.Cells.Item(31,9).Value = 1978.3043
$ws.Cells.Item(31,10).Value = 4397.6
$ws.Cells.Item(31,11).Value = 1978.3043
$ws.Cells.Item(31,12).Value = 4397.6
$ws.Cells.Item(31,13).Value = -1683.3043
$ws.Cells.Item(31,14).Value = -4987.6

# Row 34
$ws.Cells.Item(34,8).Value = 3103.558
$ws.Cells.Item(34,9).Value = 1978.3043
$ws.Cells.Item(34,10).Value = 4397.6
$ws.Cells.Item(34,11).Value = 1978.3043
$ws.Cells.Item(34,12).Value = 4397.6
$ws.Cells.Item(34,13).Value = -1776.3043
$ws.Cells.Item(34,14).Value = -4801.6

# Row 58
$ws.Cells.Item(58,8).Value = 6850998.5
$ws.Cells.Item(58,9).Value = 702.0351000000001
$ws.Cells.Item(58,10).Value = 31255180
$ws.Cells.Item(58,11).Value = 702.0351000000001
$ws.Cells.Item(58,12).Value = 31255180
$ws.Cells.Item(58,13).Value = -499.0351000000001
$ws.Cells.Item(58,14).Value = -31255586

# Row 132
$ws.Cells.Item(132,8).Value = 1623.3846
$ws.Cells.Item(132,9).Value = 1258.8857
$ws.Cells.Item(132,10).Value = 2373.8235
$ws.Cells.Item(132,11).Value = 3776.6571
$ws.Cells.Item(132,12).Value = 7121.470499999999
$ws.Cells.Item(132,13).Value = -1246.6571
$ws.Cells.Item(132,14).Value = -12181.4705

# Row 134
$ws.Cells.Item(134,8).Value = 1362.7317
$ws.Cells.Item(134,9).Value = 522.9697
$ws.Cells.Item(134,10).Value = 4826.75
$ws.Cells.Item(134,11).Value = 1568.9091
$ws.Cells.Item(134,12).Value = 14480.25
$ws.Cells.Item(134,13).Value = 966.0909000000001
$ws.Cells.Item(134,14).Value = -19550.25

# Row 136
$ws.Cells.Item(136,8).Value = 6850998.5
$ws.Cells.Item(136,9).Value = 702.0351000000001
$ws.Cells.Item(136,10).Value = 31255180
$ws.Cells.Item(136,11).Value = 2106.1053
$ws.Cells.Item(136,12).Value = 93765540
$ws.Cells.Item(136,13).Value = 443.8946999999998
$ws.Cells.Item(136,14).Value = -93770640

# Row 140
$ws.Cells.Item(140,8).Value = 44700
$ws.Cells.Item(140,9).Value = 0
$ws.Cells.Item(140,10).Value = 44700
$ws.Cells.Item(140,11).Value = 0
$ws.Cells.Item(140,12).Value = 44700
$ws.Cells.Item(140,14).Value = -55060

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131,8).Value = 1243.9183
$ws.Cells.Item(131,9).Value = 1578.8889
$ws.Cells.Item(131,10).Value = 1168.55
$ws.Cells.Item(131,11).Value = 4736.6667
$ws.Cells.Item(131,12).Value = 3505.65
$ws.Cells.Item(131,13).Value = 303.3333000000002
$ws.Cells.Item(131,14).Value = -13585.65

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80,8).Value = 2999.4
$ws.Cells.Item(80,9).Value = 2999.25
$ws.Cells.Item(80,10).Value = 3000
$ws.Cells.Item(80,11).Value = 2999.25
$ws.Cells.Item(80,12).Value = 3000
$ws.Cells.Item(80,13).Value = -2001.25
$ws.Cells.Item(80,14).Value = -4996

# Row 83
$ws.Cells.Item(83,8).Value = 2999.4
$ws.Cells.Item(83,9).Value = 2999.25
$ws.Cells.Item(83,10).Value = 3000
$ws.Cells.Item(83,11).Value = 14996.25
$ws.Cells.Item(83,12).Value = 15000
$ws.Cells.Item(83,13).Value = -10004.25
$ws.Cells.Item(83,14).Value = -24984

# Row 97
$ws.Cells.Item(97,8).Value = 2883.6316
$ws.Cells.Item(97,9).Value = 2006.5
$ws.Cells.Item(97,10).Value = 5339.6
$ws.Cells.Item(97,11).Value = 2006.5
$ws.Cells.Item(97,12).Value = 5339.6
$ws.Cells.Item(97,13).Value = -1510.5
$ws.Cells.Item(97,14).Value = -6331.6

# Row 132
$ws.Cells.Item(132,8).Value = 2134.7886
$ws.Cells.Item(132,9).Value = 1536.5625
$ws.Cells.Item(132,10).Value = 3091.95
$ws.Cells.Item(132,11).Value = 4609.6875
$ws.Cells.Item(132,12).Value = 9275.849999999999
$ws.Cells.Item(132,13).Value = -2079.6875
$ws.Cells.Item(132,14).Value = -14335.85

$ws = $wb.Worksheets.Item("LTW")
# Row 57
$ws.Cells.Item(57,8).Value = 3760.25
$ws.Cells.Item(57,9).Value = 2520.5
$ws.Cells.Item(57,10).Value = 5000
$ws.Cells.Item(57,11).Value = 2520.5
$ws.Cells.Item(57,12).Value = 5000
$ws.Cells.Item(57,13).Value = -1954.5
$ws.Cells.Item(57,14).Value = -6132

# Row 132
$ws.Cells.Item(132,8).Value = 1522.8677
$ws.Cells.Item(132,9).Value = 862.6923
$ws.Cells.Item(132,10).Value = 3668.4375
$ws.Cells.Item(132,11).Value = 2588.0769
$ws.Cells.Item(132,12).Value = 11005.3125
$ws.Cells.Item(132,13).Value = -58.07690000000002

# Row 136
$ws.Cells.Item(136,8).Value = 1632.6316
$ws.Cells.Item(136,9).Value = 1034.9333
$ws.Cells.Item(136,10).Value = 3874
$ws.Cells.Item(136,11).Value = 3104.7999
$ws.Cells.Item(136,12).Value = 11622
$ws.Cells.Item(136,13).Value = -554.7999
$ws.Cells.Item(136,14).Value = -16722

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Cells.Item(18,8).Value = 40005.25
$ws.Cells.Item(18,9).Value = 0
$ws.Cells.Item(18,10).Value = 40005.25
$ws.Cells.Item(18,11).Value = 0
$ws.Cells.Item(18,12).Value = 40005.25
$ws.Cells.Item(18,14).Value = -40351.25

# Row 132
$ws.Cells.Item(132,8).Value = 13094.957
$ws.Cells.Item(132,9).Value = 2442.8857
$ws.Cells.Item(132,10).Value = 44163.5
$ws.Cells.Item(132,11).Value = 7328.657099999999
$ws.Cells.Item(132,12).Value = 132490.5
$ws.Cells.Item(132,13).Value = -4798.657099999999
$ws.Cells.Item(132,14).Value = -137550.5
